# "Added JSE get text" - update the Doctors sheet with a new batch of
# ENT specialists (Mumbai) replacing the previous Dentist (Chennai) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

$data = @(
    @("Dr. Divya Prabhat",   "Ear-Nose-Throat (ENT) Specialist", "39 years experience overall", "Mahim,Mumbai"),
    @("Dr. Jaideep Mankani", "Ear-Nose-Throat (ENT) Specialist", "33 years experience overall", "Kandivali East,Mumbai"),
    @("Dr. Krishna Vora",    "Ear-Nose-Throat (ENT) Specialist", "30 years experience overall", "Tardeo,Mumbai"),
    @("Dr. Sonali Pandit",   "Ear-Nose-Throat (ENT) Specialist", "27 years experience overall", "Chembur,Mumbai"),
    @("Dr. Ajay Doiphode",   "Ear-Nose-Throat (ENT) Specialist", "27 years experience overall", "Andheri West,Mumbai")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
